$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Three new "blocks" of 8 rows each are appended after row 1052 (one block
# per trainee, one row per course). Each block repeats the same 8 courses
# in the same order, and the same 8 course dates for the non-numeric blocks:
#   30 Hours Construction Safety & Health   C=2    D(text)=05-01-2025
#   30 Hours G. Industry Safety & Health    C=230  D(text)=10-01-2025
#   Electrical Safety & LOTO                C=952  D(text)=06-01-2025
#   Fire Marshal                            C=262  D(text)=03-01-2025
#   Scaffold Competent Person               C=247  D(text)=01-01-2025
#   Lifting & Rigging Competent Person      C=263  D(text)=02-01-2025
#   Health & Safety Risk Assessment         C=252  D(text)=07-01-2025
#   Safety Management System & PTW          C=231  D(text)=08-01-2025
#
# Rows 869:876 and 1045:1052 are pre-existing "style twins" for the new
# blocks (same fills/borders, same C-column course order, same D-column
# text dates) - copying them keeps every cellXf index identical to the
# original rows instead of synthesizing new ones cell-by-cell.
# ---------------------------------------------------------------------------

# --- Block 1: rows 1053-1060 (Eslam Abdelsabour Khalafallah Salim) ---------
$ws.Range("A869:E876").Copy($ws.Range("A1053:E1060"))

$ws.Range("B1053").Value = "Eslam Abdelsabour Khalafallah Salim"
$ws.Range("B1054:B1060").Value = "Eslam Abdelsabour Khalafallah Salim"

$ws.Range("A1053").Value = "DSS2052"
$ws.Range("A1054").Value = "DSS2053"
$ws.Range("A1055").Value = "DSS2054"
$ws.Range("A1056").Value = "DSS2055"
$ws.Range("A1057").Value = "DSS2056"
$ws.Range("A1058").Value = "DSS2057"
$ws.Range("A1059").Value = "DSS2058"
$ws.Range("A1060").Value = "DSS2059"

# Column D here holds real date serials (new style: numFmt 14 / fill3 / border6).
# Apply the number format to a single cell first (one new cellXf), then copy
# that cell's formatting across the rest of the column so every cell shares
# the same new style instead of minting one per cell.
$ws.Range("D1053").NumberFormat = "mm-dd-yy"
$ws.Range("D1053").Copy($ws.Range("D1054:D1060"))

$ws.Range("D1053").Value = 45779
$ws.Range("D1054").Value = 45932
$ws.Range("D1055").Value = 45810
$ws.Range("D1056").Value = 45718
$ws.Range("D1057").Value = 45659
$ws.Range("D1058").Value = 45690
$ws.Range("D1059").Value = 45840
$ws.Range("D1060").Value = 45871

# --- Block 2: rows 1061-1068 (Abdullah Hesham Abdullah) --------------------
$ws.Range("A1045:E1052").Copy($ws.Range("A1061:E1068"))

$ws.Range("A1061").Value = "DSS2060"
$ws.Range("A1062").Value = "DSS2061"
$ws.Range("A1063").Value = "DSS2062"
$ws.Range("A1064").Value = "DSS2063"
$ws.Range("A1065").Value = "DSS2064"
$ws.Range("A1066").Value = "DSS2065"
$ws.Range("A1067").Value = "DSS2066"
$ws.Range("A1068").Value = "DSS2067"

$ws.Range("B1061").Value = "Abdullah Hesham Abdullah"
$ws.Range("B1062:B1068").Value = "Abdullah Hesham Abdullah"

# Column D (course dates as text, style 44) already matches verbatim via the
# copy from 1045:1052 - no further edits needed there.

# --- Block 3: rows 1069-1076 (EBNMASOUD ABDALMONIM ABDALHADI MOHAMED) ------
$ws.Range("A869:E876").Copy($ws.Range("A1069:E1076"))

$ws.Range("B1069").Value = "EBNMASOUD ABDALMONIM ABDALHADI MOHAMED"
$ws.Range("B1070:B1076").Value = "EBNMASOUD ABDALMONIM ABDALHADI MOHAMED"

$ws.Range("A1069").Value = "DSS2068"
$ws.Range("A1070").Value = "DSS2069"
$ws.Range("A1071").Value = "DSS2070"
$ws.Range("A1072").Value = "DSS2071"
$ws.Range("A1073").Value = "DSS2072"
$ws.Range("A1074").Value = "DSS2073"
$ws.Range("A1075").Value = "DSS2074"
$ws.Range("A1076").Value = "DSS2075"

# Column D (course dates as text, style 11) already matches verbatim via the
# copy from 869:876 - no further edits needed there.

# --- View state: matches the committed sheetView in the diff --------------
$ws.Range("G1076").Select()
